$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Header row rename (A1:D1) -> short machine-readable column names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2) Fix A156: drop the embedded CR/LF artifact (_x000D_ + newline) and title-case "de" -> "De"
$ws.Range("A156").Value = "Estado De México"

# Title-case Spanish connector words (de/del/la/el/los/las/y -> De/Del/La/El/Los/Las/Y)
# in municipality/state name cells.
$textUpdates = @(
    ,@('B5', 'Pabellón De Arteaga')
    ,@('B6', 'Rincón De Romos')
    ,@('B28', 'Amatenango De La Frontera')
    ,@('B35', 'Comitán De Domínguez')
    ,@('B49', 'Ocozocoautla De Espinosa')
    ,@('B54', 'Salto De Agua')
    ,@('B78', 'Guadalupe Y Calvo')
    ,@('B80', 'Hidalgo Del Parral')
    ,@('B88', 'San Francisco De Borja')
    ,@('B89', 'San Francisco Del Oro')
    ,@('B110', 'Villa De Álvarez')
    ,@('A112', 'Ciudad De México')
    ,@('B116', 'Cuajimalpa De Morelos')
    ,@('B138', 'Nombre De Dios')
    ,@('B141', 'Pánuco De Coronado')
    ,@('B145', 'San Juan De Guadalupe')
    ,@('B146', 'San Juan Del Río')
    ,@('B147', 'San Luis Del Cordero')
    ,@('B156', 'Acambay De Ruíz Castañeda')
    ,@('B159', 'Almoloya De Alquisiras')
    ,@('B160', 'Almoloya De Juárez')
    ,@('B161', 'Almoloya Del Río')
    ,@('B167', 'Atizapán De Zaragoza')
    ,@('B172', 'Chapa De Mota')
    ,@('B181', 'Ecatepec De Morelos')
    ,@('B186', 'Ixtapan De La Sal')
    ,@('B197', 'Naucalpan De Juárez')
    ,@('B206', 'San Antonio La Isla')
    ,@('B207', 'San Felipe Del Progreso')
    ,@('B208', 'San Martín De Las Pirámides')
    ,@('B209', 'San Simón De Guerrero')
    ,@('B210', 'Soyaniquilpan De Juárez')
    ,@('B217', 'Tenango Del Valle')
    ,@('B226', 'Tlalnepantla De Baz')
    ,@('B232', 'Valle De Bravo')
    ,@('B233', 'Valle De Chalco Solidaridad')
    ,@('B234', 'Villa De Allende')
    ,@('B235', 'Villa Del Carbón')
    ,@('B244', 'Apaseo El Alto')
    ,@('B252', 'Dolores Hidalgo Cuna De La Independencia Nacional')
    ,@('B256', 'Jaral Del Progreso')
    ,@('B263', 'Purísima Del Rincón')
    ,@('B267', 'San Diego De La Unión')
    ,@('B269', 'San Francisco Del Rincón')
    ,@('B271', 'San Luis De La Paz')
    ,@('B272', 'Santa Cruz De Juventino Rosas')
    ,@('B274', 'Silao De La Victoria')
    ,@('B278', 'Valle De Santiago')
    ,@('B283', 'Acapulco De Juárez')
    ,@('B286', 'Ajuchitlán Del Progreso')
    ,@('B287', 'Alcozauca De Guerrero')
    ,@('B290', 'Atenango Del Río')
    ,@('B292', 'Atoyac De Álvarez')
    ,@('B293', 'Ayutla De Los Libres')
    ,@('B295', 'Buenavista De Cuéllar')
    ,@('B296', 'Chilapa De Álvarez')
    ,@('B297', 'Chilpancingo De Los Bravo')
    ,@('B298', 'Coahuayutla De José María Izazaga')
    ,@('B303', 'Coyuca De Benítez')
    ,@('B304', 'Coyuca De Catalán')
    ,@('B308', 'Cuetzala Del Progreso')
    ,@('B309', 'Cutzamala De Pinzón')
    ,@('B315', 'Huitzuco De Los Figueroa')
    ,@('B316', 'Iguala De La Independencia')
    ,@('B318', 'Ixcateopan De Cuauhtémoc')
    ,@('B321', 'La Unión De Isidoro Montes De Oca')
    ,@('B324', 'Mártir De Cuilapan')
    ,@('B336', 'Taxco De Alarcón')
    ,@('B338', 'Técpan De Galeana')
    ,@('B340', 'Tepecoacuilco De Trujano')
    ,@('B342', 'Tixtla De Guerrero')
    ,@('B345', 'Tlalixtaquilla De Maldonado')
    ,@('B346', 'Tlapa De Comonfort')
    ,@('B357', 'Agua Blanca De Iturbide')
    ,@('B361', 'Atotonilco El Grande')
    ,@('B366', 'Cuautepec De Hinojosa')
    ,@('B371', 'Huasca De Ocampo')
    ,@('B374', 'Huejutla De Reyes')
    ,@('B377', 'Jacala De Ledezma')
    ,@('B382', 'Mineral Del Chico')
    ,@('B383', 'Mineral Del Monte')
    ,@('B384', 'Mixquiahuala De Juárez')
    ,@('B386', 'Omitlán De Juárez')
    ,@('B387', 'Pachuca De Soto')
    ,@('B388', 'Progreso De Obregón')
    ,@('B393', 'Santiago De Anaya')
    ,@('B394', 'Santiago Tulantepec De Lugo Guerrero')
    ,@('B397', 'Tenango De Doria')
    ,@('B399', 'Tepehuacán De Guerrero')
    ,@('B402', 'Tezontepec De Aldama')
    ,@('B407', 'Tula De Allende')
    ,@('B408', 'Tulancingo De Bravo')
    ,@('B410', 'Zacualtipán De Ángeles')
    ,@('B411', 'Zapotlán De Juárez')
    ,@('B415', 'Acatlán De Juárez')
    ,@('B416', 'Ahualulco De Mercado')
    ,@('B421', 'Atemajac De Brizuela')
    ,@('B424', 'Atotonilco El Alto')
    ,@('B426', 'Autlán De Navarro')
    ,@('B438', 'Concepción De Buenos Aires')
    ,@('B439', 'Cuautitlán De García Barragán')
    ,@('B446', 'Encarnación De Díaz')
    ,@('B453', 'Huejuquilla El Alto')
    ,@('B454', 'Ixtlahuacán De Los Membrillos')
    ,@('B455', 'Ixtlahuacán Del Río')
    ,@('B459', 'Jilotlán De Los Dolores')
    ,@('B464', 'Lagos De Moreno')
    ,@('B472', 'Ojuelos De Jalisco')
    ,@('B478', 'San Juan De Los Lagos')
    ,@('B480', 'San Martín De Bolaños')
    ,@('B482', 'San Miguel El Alto')
    ,@('B483', 'San Sebastián Del Oeste')
    ,@('B484', 'Santa María De Los Ángeles')
    ,@('B487', 'Talpa De Allende')
    ,@('B488', 'Tamazula De Gordiano')
    ,@('B494', 'Teocuitatlán De Corona')
    ,@('B495', 'Tepatitlán De Morelos')
    ,@('B498', 'Tizapán El Alto')
    ,@('B499', 'Tlajomulco De Zúñiga')
    ,@('B510', 'Unión De Tula')
    ,@('B511', 'Valle De Guadalupe')
    ,@('B515', 'Yahualica De González Gallo')
    ,@('B516', 'Zacoalco De Torres')
    ,@('B519', 'Zapotitlán De Vadillo')
    ,@('B520', 'Zapotlán Del Rey')
    ,@('B521', 'Zapotlán El Grande')
    ,@('B545', 'Coalcomán De Vázquez Pallares')
    ,@('B547', 'Cojumatlán De Régules')
    ,@('B613', 'Tiquicheo De Nicolás Romero')
    ,@('B637', 'Coatlán Del Río')
    ,@('B647', 'Puente De Ixtla')
    ,@('B651', 'Tetela Del Volcán')
    ,@('B652', 'Tlaltizapán De Zapata')
    ,@('B662', 'Amatlán De Cañas')
    ,@('B663', 'Bahía De Banderas')
    ,@('B667', 'Ixtlán Del Río')
    ,@('B674', 'Santa María Del Oro')
    ,@('B687', 'San Nicolás De Los Garza')
    ,@('B691', 'Acatlán De Pérez Figueroa')
    ,@('B698', 'Chalcatongo De Hidalgo')
    ,@('B699', 'Ciénega De Zimatlán')
    ,@('B702', 'Coicoyán De Las Flores')
    ,@('B703', 'Constancia Del Rosario')
    ,@('B705', 'Cuilápam De Guerrero')
    ,@('B706', 'Eloxochitlán De Flores Magón')
    ,@('B707', 'Guadalupe De Ramírez')
    ,@('B709', 'Heroica Ciudad De Ejutla De Crespo')
    ,@('B710', 'Heroica Ciudad De Huajuapan De León')
    ,@('B711', 'Heroica Ciudad De Tlaxiaco')
    ,@('B713', 'Ixtlán De Juárez')
    ,@('B714', 'Heroica Ciudad De Juchitán De Zaragoza')
    ,@('B723', 'Mariscala De Juárez')
    ,@('B726', 'Miahuatlán De Porfirio Díaz')
    ,@('B728', 'Oaxaca De Juárez')
    ,@('B729', 'Ocotlán De Morelos')
    ,@('B730', 'Pinotepa De Don Luis')
    ,@('B732', 'Putla Villa De Guerrero')
    ,@('B779', 'San Juan Del Estado')
    ,@('B780', 'San Juan Del Río')
    ,@('B813', 'San Miguel El Grande')
    ,@('B826', 'San Pablo Villa De Mitla')
    ,@('B831', 'San Pedro El Alto')
    ,@('B844', 'San Pedro Y San Pablo Ayutla')
    ,@('B853', 'Santa Ana Del Valle')
    ,@('B862', 'Santa Cruz Tacache De Mina')
    ,@('B867', 'Santa Inés Del Monte')
    ,@('B877', 'Santa María Jalapa Del Marqués')
    ,@('B917', 'Santo Domingo De Morelos')
    ,@('B932', 'Tamazulápam Del Espíritu Santo')
    ,@('B933', 'Tataltepec De Valdés')
    ,@('B934', 'Teococuilco De Marcos Pérez')
    ,@('B935', 'Teotitlán De Flores Magón')
    ,@('B936', 'Teotitlán Del Valle')
    ,@('B937', 'Tepelmeme Villa De Morelos')
    ,@('B939', 'Tlacolula De Matamoros')
    ,@('B940', 'Totontepec Villa De Morelos')
    ,@('B942', 'Villa De Chilapa De Díaz')
    ,@('B943', 'Villa De Etla')
    ,@('B944', 'Villa De Tamazulápam Del Progreso')
    ,@('B945', 'Villa De Tututepec De Melchor Ocampo')
    ,@('B946', 'Villa De Zaachila')
    ,@('B949', 'Villa Sola De Vega')
    ,@('B955', 'Zimatlán De Álvarez')
    ,@('B987', 'Cuayuca De Andrade')
    ,@('B997', 'Huehuetlán El Chico')
    ,@('B1002', 'Ixcamilpa De Guerrero')
    ,@('B1005', 'Izúcar De Matamoros')
    ,@('B1010', 'Los Reyes De Juárez')
    ,@('B1031', 'San Salvador El Seco')
    ,@('B1034', 'Tecali De Herrera')
    ,@('B1042', 'Tepanco De López')
    ,@('B1043', 'Tepatlaxco De Hidalgo')
    ,@('B1045', 'Tepexi De Rodríguez')
    ,@('B1047', 'Tepeyahualco De Cuauhtémoc')
    ,@('B1048', 'Tetela De Ocampo')
    ,@('B1052', 'Tlacotepec De Benito Juárez')
    ,@('B1061', 'Totoltepec De Guerrero')
    ,@('B1066', 'Xayacatlán De Bravo')
    ,@('B1070', 'Xochitlán De Vicente Suárez')
    ,@('B1080', 'Amealco De Bonfil')
    ,@('B1081', 'Cadereyta De Montes')
    ,@('B1085', 'Jalpan De Serra')
    ,@('B1086', 'Landa De Matamoros')
    ,@('B1089', 'Pinal De Amoles')
    ,@('B1092', 'San Juan Del Río')
    ,@('B1109', 'Mexquitic De Carmona')
    ,@('B1113', 'San Ciro De Acosta')
    ,@('B1115', 'Santa María Del Río')
    ,@('B1117', 'Soledad De Graciano Sánchez')
    ,@('B1120', 'Villa De Ramos')
    ,@('B1158', 'Nacozari De García')
    ,@('B1164', 'San Pedro De La Cueva')
    ,@('B1189', 'Soto La Marina')
    ,@('B1192', 'Acuamanala De Miguel Hidalgo')
    ,@('B1193', 'Apetatitlán De Antonio Carvajal')
    ,@('B1200', 'Ixtacuixtla De Mariano Matamoros')
    ,@('B1205', 'San Pablo Del Monte')
    ,@('B1206', 'Tetla De La Solidaridad')
    ,@('B1222', 'Alto Lucero De Gutiérrez Barrios')
    ,@('B1225', 'Amatlán De Los Reyes')
    ,@('B1233', 'Boca Del Río')
    ,@('B1238', 'Castillo De Teayo')
    ,@('B1240', 'Cazones De Herrera')
    ,@('B1249', 'Cosamaloapan De Carpio')
    ,@('B1263', 'Hueyapan De Ocampo')
    ,@('B1264', 'Ignacio De La Llave')
    ,@('B1267', 'Ixhuatlán De Madero')
    ,@('B1268', 'Ixhuatlán Del Café')
    ,@('B1275', 'Juchique De Ferrer')
    ,@('B1279', 'Lerdo De Tejada')
    ,@('B1283', 'Martínez De La Torre')
    ,@('B1284', 'Medellín De Bravo')
    ,@('B1297', 'Paso Del Macho')
    ,@('B1300', 'Poza Rica De Hidalgo')
    ,@('B1306', 'Sayula De Alemán')
    ,@('B1309', 'Soledad De Doblado')
    ,@('B1337', 'Zontecomatlán De López Y Fuentes')
    ,@('B1338', 'Zozocolco De Hidalgo')
    ,@('B1370', 'Concepción Del Oro')
    ,@('B1372', 'El Plateado De Joaquín Amaro')
    ,@('B1380', 'Jiménez Del Teul')
    ,@('B1388', 'Moyahua De Estrada')
    ,@('B1389', 'Nochistlán De Mejía')
    ,@('B1390', 'Noria De Ángeles')
    ,@('B1399', 'Teúl De González Ortega')
    ,@('B1400', 'Tlaltenango De Sánchez Román')
    ,@('B1401', 'Trinidad García De La Cadena')
    ,@('B1404', 'Villa De Cos')
)
foreach ($u in $textUpdates) {
    $ws.Range($u[0]).Value = $u[1]
}

# D-column micro precision fix (last-ULP float adjustments carried over from the source file)
$dUpdates = @(
    ,@('D17', 0.0009447628645210052)
    ,@('D58', 0.0009447628645210052)
    ,@('D207', 0.0009447628645210052)
    ,@('D228', 0.0009447628645210052)
    ,@('D291', 0.0009447628645210052)
    ,@('D298', 0.0009447628645210052)
    ,@('D331', 0.0009447628645210052)
    ,@('D336', 0.0009447628645210052)
    ,@('D351', 0.0009447628645210052)
    ,@('D380', 0.0009447628645210052)
    ,@('D545', 0.0009447628645210052)
    ,@('D619', 0.0009447628645210052)
    ,@('D676', 0.0009447628645210052)
    ,@('D686', 0.0009447628645210052)
    ,@('D691', 0.0009447628645210052)
    ,@('D837', 0.0009447628645210052)
    ,@('D924', 0.0009447628645210052)
    ,@('D968', 0.0009447628645210052)
    ,@('D1005', 0.0009447628645210052)
    ,@('D1135', 0.0009447628645210052)
    ,@('D1194', 0.0009447628645210052)
    ,@('D1317', 0.0009447628645210052)
    ,@('D1216', 0.009258676072305853)
)
foreach ($u in $dUpdates) {
    $ws.Range($u[0]).Value = $u[1]
}

# 3) Delete the trailing footnote rows (1411-1415); row 1410 is already blank,
#    so the used range collapses to A1:D1409 automatically.
$ws.Range("A1411:D1415").EntireRow.Delete()
